$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.707.77'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.600.02'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.31'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("E6").Value = '  -0.37%  '

$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("E8").Value = '  +0.37%  '

$ws.Range("E9").Value = '  +0.98%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.53'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("D12").Value = '1.824.63'
$ws.Range("E12").Value = '  +0.40%  '

$ws.Range("D13").Value = '1.607.44'
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("E14").Value = '  +0.59%  '

$ws.Range("E15").Value = '  +0.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.35'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.44%  '

$ws.Range("D17").Value = '26.687.28'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").Value = '0.0₃0755'
$ws.Range("E18").Value = '  +3.63%  '

$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.43'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.18'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.67%  '

$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.92'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("E27").Value = '  +0.48%  '

$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("E29").Value = '  +0.64%  '

$ws.Range("E30").Value = '  +2.66%  '

$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.78%  '

$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("D34").Value = '1.291.20'
$ws.Range("E34").Value = '  +1.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.619'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.12%  '

$ws.Range("E36").Value = '  +1.07%  '

$ws.Range("E37").Value = '  +0.52%  '

$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("E39").Value = '  +16.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.825'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.75%  '

$ws.Range("E41").Value = '  -0.90%  '

$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("E43").Value = '  -0.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.10'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.36%  '

$ws.Range("D45").Value = '1.735.79'
$ws.Range("E45").Value = '  +0.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.94'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.59%  '

$ws.Range("E47").Value = '  -0.97%  '

$ws.Range("E48").Value = '  -1.00%  '

$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("E50").Value = '  +0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.35'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.12%  '
